$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.100.45"
$ws.Range("E2").Value = "  +4.17%  "
$ws.Range("D3").Value = "2.635.69"
$ws.Range("E3").Value = "  +5.47%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'327.90"
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("D6").Value = "'110.72"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.561"
$ws.Range("E9").Value = "  +4.22%  "
$ws.Range("D10").Value = "'40.53"
$ws.Range("E10").Value = "  +2.32%  "
$ws.Range("D11").Value = "'20.71"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("D12").Value = "'0.0821"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "'7.30"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "3.050.81"
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("D16").Value = "2.642.59"
$ws.Range("E16").Value = "  +5.93%  "
$ws.Range("D17").Value = "'0.878"
$ws.Range("E17").Value = "  +4.94%  "
$ws.Range("D18").Value = "49.988.80"
$ws.Range("E18").Value = "  +4.29%  "
$ws.Range("D19").Value = "'3.09"
$ws.Range("E19").Value = "  +13.02%  "
$ws.Range("D20").Value = "'13.33"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("D21").Value = "'6.83"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "'72.97"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "'279.98"
$ws.Range("E24").Value = "  +1.02%  "
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "'26.64"
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'36.87"
$ws.Range("E28").Value = "  +5.33%  "
$ws.Range("E29").Value = "  +5.91%  "
$ws.Range("E30").Value = "  +1.98%  "
$ws.Range("D31").Value = "'0.144"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").Value = "'49.86"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'19.82"
$ws.Range("E33").Value = "  +1.04%  "
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").Value = "'2.06"
$ws.Range("E37").Value = "  +6.31%  "
$ws.Range("D38").Value = "'4.75"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("E39").Value = "  +7.22%  "
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Value = "'123.13"
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("D42").Value = "'22.46"
$ws.Range("E42").Value = "  +5.08%  "
$ws.Range("D43").Value = "'2.22"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("E44").Value = "  +4.47%  "
$ws.Range("E45").Value = "  +6.54%  "
$ws.Range("D46").Value = "2.059.91"
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("D47").Value = "'2.34"
$ws.Range("E47").Value = "  +17.51%  "
$ws.Range("D48").Value = "'2.02"
$ws.Range("E48").Value = "  +9.22%  "
$ws.Range("D49").Value = "'9.06"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "'5.38"
$ws.Range("E50").Value = "  +4.15%  "
$ws.Range("D51").Value = "'81.88"
$ws.Range("E51").Value = "  +2.10%  "
